$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: split "Primeiro Método, Cliente ocupado, com ele reunião de tempo
# limitado" into three runs by inserting a new word "fazer " between
# "Cliente ocupado, " and "com ele reunião de tempo limitado".
# ---------------------------------------------------------------------------
$f1 = $d.Content
$f1.Find.Execute("com ele reunião de tempo limitado") | Out-Null
$tailStart = $f1.Start
$tailEnd = $f1.End

# Clear the tail text, then re-insert "fazer " and the tail as two separate
# InsertAfter calls so each becomes its own run (no leftover rPr).
$tailRange = $d.Range($tailStart, $tailEnd)
$tailRange.Text = ""

$p2 = $d.Range($tailStart, $tailStart)
$p2.InsertAfter("fazer ")

$p3 = $d.Range($p2.End, $p2.End)
$p3.InsertAfter("com ele reunião de tempo limitado")

# ---------------------------------------------------------------------------
# Edit 2: merge the four runs of "2 - Quem são os maiores concorrentes..."
# into a single run by replacing the whole sentence via Find/Replace.
# ---------------------------------------------------------------------------
$f2 = $d.Content
$f2.Find.Execute(
    "2 - Quem são os maiores concorrentes e o que os preocupa sobre eles? - Analise de competição",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2 - Quem são os maiores concorrentes e o que os preocupa sobre eles? - Analise de competição",
    2) | Out-Null

# ---------------------------------------------------------------------------
# Edit 3: move the "_GoBack" bookmark from its own paragraph near the end of
# the document to right after "Cliente tem um dia todo disponível".
# ---------------------------------------------------------------------------
$f3 = $d.Content
$f3.Find.Execute("Cliente tem um dia todo disponível") | Out-Null
$tailPos = $f3.End

# Adding a bookmark collapsed exactly on a paragraph boundary misbehaves, so
# insert a temporary marker character after the text to give the insertion
# point a safe (non-boundary) position, add the bookmark there, then remove
# the marker again.
$marker = $d.Range($tailPos, $tailPos)
$marker.InsertAfter("X")

$bmTarget = $d.Range($tailPos, $tailPos)
$d.Bookmarks.Add("_GoBack", $bmTarget) | Out-Null

$bm = $d.Bookmarks.Item("_GoBack")
$markerRange = $d.Range($bm.End, $bm.End + 1)
$markerRange.Text = ""
